$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1 ("time_taken"), matching the bold/centered/bordered
# header style already used by B1:E1 (copy E1's formatting onto F1).
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# New data cell F2 holding the recorded timestamp as plain text.
$ws.Range("F2").Value = "2021-10-05 13:38:33.327531"
